# Natmi LR-pair table (Clcf1 -> Il6st) following Dr Hou's advice:
# recompute the Sending-cluster x Target-cluster grid for ECs / FAPs / sCs
# (3x3 = 9 data rows instead of the previous 3x2 = 6 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Clcf1 -> Il6st -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Clcf1"
$ws.Range("C2").Value = "Il6st"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.9214586666666666
$ws.Range("H2").Value = 2.764376
$ws.Range("I2").Value = 0.08041853843186561
$ws.Range("J2").Value = 0.08041853843186561
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 39.96608766666667
$ws.Range("N2").Value = 119.898263
$ws.Range("O2").Value = 0.2616165719423124
$ws.Range("P2").Value = 0.2616165719423124
$ws.Range("Q2").Value = 36.82709785320979
$ws.Range("R2").Value = 331.443880678888
$ws.Range("S2").Value = 0.02103882234515578
$ws.Range("T2").Value = 0.02103882234515578

# Row 3: ECs -> Clcf1 -> Il6st -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Clcf1"
$ws.Range("C3").Value = "Il6st"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.9214586666666666
$ws.Range("H3").Value = 2.764376
$ws.Range("I3").Value = 0.08041853843186561
$ws.Range("J3").Value = 0.08041853843186561
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 85.11961100000001
$ws.Range("N3").Value = 255.358833
$ws.Range("O3").Value = 0.5571899111219771
$ws.Range("P3").Value = 0.557189911121977
$ws.Range("Q3").Value = 78.43420325924534
$ws.Range("R3").Value = 705.907829333208
$ws.Range("S3").Value = 0.0448083982814105
$ws.Range("T3").Value = 0.04480839828141049

# Row 4: ECs -> Clcf1 -> Il6st -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Clcf1"
$ws.Range("C4").Value = "Il6st"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.9214586666666666
$ws.Range("H4").Value = 2.764376
$ws.Range("I4").Value = 0.08041853843186561
$ws.Range("J4").Value = 0.08041853843186561
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 27.68018833333333
$ws.Range("N4").Value = 83.040565
$ws.Range("O4").Value = 0.1811935169357105
$ws.Range("P4").Value = 0.1811935169357105
$ws.Range("Q4").Value = 25.50614943471556
$ws.Range("R4").Value = 229.55534491244
$ws.Range("S4").Value = 0.01457131780529932
$ws.Range("T4").Value = 0.01457131780529932

# Row 5: FAPs -> Clcf1 -> Il6st -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Clcf1"
$ws.Range("C5").Value = "Il6st"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.913147
$ws.Range("H5").Value = 5.739441
$ws.Range("I5").Value = 0.1669662363715809
$ws.Range("J5").Value = 0.1669662363715809
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 39.96608766666667
$ws.Range("N5").Value = 119.898263
$ws.Range("O5").Value = 0.2616165719423124
$ws.Range("P5").Value = 0.2616165719423124
$ws.Range("Q5").Value = 76.46100072122036
$ws.Range("R5").Value = 688.1490064909831
$ws.Range("S5").Value = 0.04368113438964283
$ws.Range("T5").Value = 0.04368113438964283

# Row 6: FAPs -> Clcf1 -> Il6st -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Clcf1"
$ws.Range("C6").Value = "Il6st"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.913147
$ws.Range("H6").Value = 5.739441
$ws.Range("I6").Value = 0.1669662363715809
$ws.Range("J6").Value = 0.1669662363715809
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 85.11961100000001
$ws.Range("N6").Value = 255.358833
$ws.Range("O6").Value = 0.5571899111219771
$ws.Range("P6").Value = 0.557189911121977
$ws.Range("Q6").Value = 162.846328425817
$ws.Range("R6").Value = 1465.616955832353
$ws.Range("S6").Value = 0.09303190240425217
$ws.Range("T6").Value = 0.09303190240425215

# Row 7: FAPs -> Clcf1 -> Il6st -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Clcf1"
$ws.Range("C7").Value = "Il6st"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.913147
$ws.Range("H7").Value = 5.739441
$ws.Range("I7").Value = 0.1669662363715809
$ws.Range("J7").Value = 0.1669662363715809
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.68018833333333
$ws.Range("N7").Value = 83.040565
$ws.Range("O7").Value = 0.1811935169357105
$ws.Range("P7").Value = 0.1811935169357105
$ws.Range("Q7").Value = 52.95626926935167
$ws.Range("R7").Value = 476.606423424165
$ws.Range("S7").Value = 0.03025319957768587
$ws.Range("T7").Value = 0.03025319957768587

# Row 8 (new): sCs -> Clcf1 -> Il6st -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Clcf1"
$ws.Range("C8").Value = "Il6st"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.623680999999999
$ws.Range("H8").Value = 25.871043
$ws.Range("I8").Value = 0.7526152251965536
$ws.Range("J8").Value = 0.7526152251965536
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 39.96608766666667
$ws.Range("N8").Value = 119.898263
$ws.Range("O8").Value = 0.2616165719423124
$ws.Range("P8").Value = 0.2616165719423124
$ws.Range("Q8").Value = 344.6547908553677
$ws.Range("R8").Value = 3101.893117698309
$ws.Range("S8").Value = 0.1968966152075138
$ws.Range("T8").Value = 0.1968966152075138

# Row 9 (new): sCs -> Clcf1 -> Il6st -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Clcf1"
$ws.Range("C9").Value = "Il6st"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.623680999999999
$ws.Range("H9").Value = 25.871043
$ws.Range("I9").Value = 0.7526152251965536
$ws.Range("J9").Value = 0.7526152251965536
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 85.11961100000001
$ws.Range("N9").Value = 255.358833
$ws.Range("O9").Value = 0.5571899111219771
$ws.Range("P9").Value = 0.557189911121977
$ws.Range("Q9").Value = 734.044372108091
$ws.Range("R9").Value = 6606.399348972819
$ws.Range("S9").Value = 0.4193496104363145
$ws.Range("T9").Value = 0.4193496104363144

# Row 10 (new): sCs -> Clcf1 -> Il6st -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Clcf1"
$ws.Range("C10").Value = "Il6st"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.623680999999999
$ws.Range("H10").Value = 25.871043
$ws.Range("I10").Value = 0.7526152251965536
$ws.Range("J10").Value = 0.7526152251965536
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 27.68018833333333
$ws.Range("N10").Value = 83.040565
$ws.Range("O10").Value = 0.1811935169357105
$ws.Range("P10").Value = 0.1811935169357105
$ws.Range("Q10").Value = 238.7051142065883
$ws.Range("R10").Value = 2148.346027859295
$ws.Range("S10").Value = 0.1363689995527253
$ws.Range("T10").Value = 0.1363689995527253
